$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row (row 1) - new columns I and J, using same style as existing header cells (s="1")
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Data values for columns I (I0) and J (IF), rows 2-16
$dataI = @(8, 9, 7, 5, 7, 7, 7, 6, 8, 9, 11, 7, 7, 9, 9)
$dataJ = @(9, 9, 7, 6, 7, 8, 8, 7, 8, 9, 11, 8, 7, 9, 9)

for ($i = 0; $i -lt $dataI.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 9).Value = $dataI[$i]
    $ws.Cells.Item($row, 10).Value = $dataJ[$i]
}
